# Fill in the monthly figures for each year row (2015-2022) on Sheet1.
# B:M = JAN..DEC, N = Total. Rows already carry their number formatting/
# borders (styles s="2"/"s=3"/"s=7" etc.), so only the values need setting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(53661, 94274, 9719,  9561,  43220, 70490, 13207, 51277, 17749, 81497, 81910, 7637,  51365)  # 2015
    3 = @(97609, 82806, 95048, 91711, 44471, 38191, 51861, 7666,  68885, 28398, 6291,  76229, 19020)  # 2016
    4 = @(92377, 93875, 10848, 28188, 14053, 68642, 51904, 36298, 10080, 84897, 65015, 63808, 42805)  # 2017
    5 = @(19343, 23299, 12058, 74382, 94679, 57963, 62805, 8784,  56526, 83635, 56246, 28584, 22366)  # 2018
    6 = @(94812, 92439, 73680, 71303, 92339, 46050, 28690, 24254, 63483, 35015, 54918, 98967, 23997)  # 2019
    7 = @(33801, 29687, 97658, 95271, 85309, 50297, 30784, 28751, 31791, 16397, 88821, 73233, 52286)  # 2020
    8 = @(82034, 67626, 59222, 63627, 94583, 94660, 93544, 13009, 87646, 13428, 51805, 22277, 21258)  # 2021
    9 = @(87142, 51384, 11236, 63681, 75474, 26524, 96017, 65310, 32117, 19081, 34257, 37292, 27952)  # 2022
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2   # column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Move the active selection from the old R9 to N10, as it was left after entry.
$ws.Range("N10").Select() | Out-Null
